$d = $word.ActiveDocument

function Insert-ParaXml($range, $innerXml) {
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# Paragraph 14 ("B. Each solution is different for every case.") is the very
# last paragraph in the body. It keeps its ListParagraph style but loses its
# run entirely (the sentence moves up into the new paragraph inserted after
# paragraph 13). Trim the paragraph mark off the range before clearing the
# text so we only wipe the run, not the body's terminating mark.
# ---------------------------------------------------------------------------
$p14 = $d.Paragraphs(14).Range
$p14.End = $p14.End - 1
$p14.Text = ""

# ---------------------------------------------------------------------------
# Paragraph 13 (numbered "A. Each solution does seem to meet the goals.")
# expands into three new paragraphs:
#   - "4. " + "A. Each solution does seem to meet the goals."
#   - "B. Each solution is different for every case." (ListParagraph, no numbering)
#   - "5. A. the solution I choose is..." + the _GoBack bookmark
# ---------------------------------------------------------------------------
$p13 = $d.Paragraphs(13).Range
$inner13 = '<w:p><w:r><w:t xml:space="preserve">      4. </w:t></w:r><w:r><w:t>A. Each solution does seem to meet the goals.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>B. Each solution is different for every case.</w:t></w:r></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">      5. A. the solution I choose is to leave the cat and the seed on shore while the man      travels with the bird across the river.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
Insert-ParaXml $p13 $inner13

# ---------------------------------------------------------------------------
# Paragraph 10 ("       3. A. Cat and seed- seems like correct choice") gets
# its leading-whitespace run split in two: a run with six spaces, then a
# separate run holding "3. A." (previously one run held all seven spaces
# plus "3. A.").
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs(10).Range
$inner10 = '<w:p><w:r><w:t xml:space="preserve">      </w:t></w:r><w:r><w:t>3. A.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Cat and seed- seems like correct choice</w:t></w:r></w:p>'
Insert-ParaXml $p10 $inner10

Write-Output "ok"
